# Applies the 27-Jun-2020 21:22 COVID data refresh to the "Pais" sheet.
# The feed re-sorts countries by "Casos totales" (column B) descending, so a
# handful of countries that overtook their neighbours change row, and all the
# refreshed totals/actives/recovered/deaths figures are written in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# "Datos actualizados a ..." banner in A1
$ws.Range("A1").Value2 = "Datos actualizados a 27 de Junio de 2020 a las 21:22"

# Row => [Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes]
$updatedRows = @(
    @{ Row = 4; Country = "Estados Unidos"; Stats = @(2582008, 29052, 1071438, 1382597, 0, 333, 127973) }
    @{ Row = 7; Country = "India"; Stats = @(529577, 20131, 310146, 203328, 0, 414, 16103) }
    @{ Row = 31; Country = "Ecuador"; Stats = @(53856, 0, 26493, 22939, 0, 18, 4424) }
    @{ Row = 52; Country = "Israel"; Stats = @(23421, 621, 17002, 6102, 0, 3, 317) }
    @{ Row = 53; Country = "Nigeria"; Stats = @(23298, 0, 8253, 14491, 0, 0, 554) }
    @{ Row = 74; Country = "Uzbekistan"; Stats = @(7682, 255, 5240, 2422, 0, 0, 20) }
    @{ Row = 75; Country = "Australia"; Stats = @(7641, 46, 6979, 558, 0, 0, 104) }
    @{ Row = 99; Country = "Costa Rica"; Stats = @(2979, 143, 1325, 1641, 0, 1, 13) }
    @{ Row = 100; Country = "Somalia"; Stats = @(2878, 0, 868, 1920, 0, 0, 90) }
    @{ Row = 114; Country = "Estado de Palestina"; Stats = @(1815, 258, 446, 1365, 0, 1, 4) }
    @{ Row = 125; Country = "Tunez"; Stats = @(1168, 4, 1025, 93, 0, 0, 50) }
    @{ Row = 129; Country = "Cabo Verde"; Stats = @(1091, 64, 568, 511, 0, 2, 12) }
    @{ Row = 130; Country = "Congo"; Stats = @(1087, 0, 456, 594, 0, 0, 37) }
    @{ Row = 131; Country = "Niger"; Stats = @(1062, 3, 924, 71, 0, 0, 67) }
    @{ Row = 132; Country = "Benin"; Stats = @(1053, 0, 292, 747, 0, 0, 14) }
    @{ Row = 143; Country = "Suazilandia"; Stats = @(745, 17, 370, 367, 0, 0, 8) }
    @{ Row = 144; Country = "Liberia"; Stats = @(729, 45, 291, 404, 0, 0, 34) }
    @{ Row = 177; Country = "Namibia"; Stats = @(136, 15, 22, 114, 0, 0, 0) }
    @{ Row = 178; Country = "Trinidad yTobago"; Stats = @(124, 0, 109, 7, 0, 0, 8) }
    @{ Row = 201; Country = "Laos"; Stats = @(19, 0, 19, 0, 0, 0, 0) }
    @{ Row = 202; Country = "Santa Lucia"; Stats = @(19, 0, 19, 0, 0, 0, 0) }
    @{ Row = 212; Country = "Seychelles"; Stats = @(11, 0, 11, 0, 0, 0, 0) }
    @{ Row = 213; Country = "Montserrat"; Stats = @(11, 0, 10, 0, 0, 0, 1) }
)

foreach ($entry in $updatedRows) {
    $ws.Cells.Item($entry.Row, 1).Value2 = $entry.Country
    for ($col = 0; $col -lt $entry.Stats.Length; $col++) {
        $ws.Cells.Item($entry.Row, $col + 2).Value2 = $entry.Stats[$col]
    }
}
